# Tutorial blocks were added: insert two "practice" condition rows
# (PID = -1) at the top of the data table, then re-sort the whole
# table by column A (PID), which is how the original author's
# spreadsheet ended up with the practice rows first followed by the
# original R1..R4 rows (now shifted down two rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new blank rows right after the header row (rows 2 & 3),
# pushing the existing data down.
$ws.Rows("2:3").Insert() | Out-Null

# Fill in the new practice rows.
$ws.Range("A2").Value = -1
$ws.Range("B2").Value = "practice1"
$ws.Range("C2").Value = "ND"

$ws.Range("A3").Value = -1
$ws.Range("B3").Value = "practice2"
$ws.Range("C3").Value = "D"

# Re-sort the data range (A1:C8, header included) by column A ascending,
# exactly as if the user selected the table and used Data > Sort.
$sortRange = $ws.Range("A1:C8")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A1:A8"), 0, 1, 0, 0) | Out-Null
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 1
$ws.Sort.MatchCase = $false
$ws.Sort.Orientation = 1
$ws.Sort.SortMethod = 1
$ws.Sort.Apply() | Out-Null

# Leave the selection where the author's session ended up.
$ws.Range("P8").Select() | Out-Null
